# Apply the "Add files via upload" revision to the Project Expertise workbook.
#
# Substance of the change (the only parts that are deterministically
# reproducible via the Excel object model -- GUIDs / revision ids / sheet
# protection salts are regenerated by Excel itself on every save and are not
# something a COM script can control):
#   1. Sheet "9 Gerard Cutright" (2nd worksheet): update several "Expertise"
#      ratings in column G (rows 2-26) - some cleared, some changed to a
#      different letter, one new "H" rating introduced (which becomes a new
#      shared string).
#   2. The active selection on that sheet moves from D23 to A10:K10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("9 Gerard Cutright")

# --- Column G "Expertise" rating updates (row -> new value) ---------------
$ratings = @{
    4  = ""
    5  = "L"
    7  = "H"
    9  = ""
    13 = "M"
    14 = "L"
    15 = "L"
    16 = "M"
    17 = "M"
    18 = "L"
    19 = "L"
    21 = "H"
    23 = "L"
    24 = "L"
    25 = "L"
    26 = ""
}

foreach ($row in $ratings.Keys) {
    $ws.Range("G$row").Value2 = $ratings[$row]
}

# --- Update the (unfrozen top pane) selection ------------------------------
$ws.Activate()
$ws.Range("A10:K10").Select()
